$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. All target cells hold text-formatted
# numbers/percentages (e.g. "274.87", "0.46%"), so we force the
# NumberFormat to Text before assignment (otherwise Excel auto-converts
# numeric-looking strings like "274.87" into real numbers), then restore
# the "Normal" style afterwards so the cell keeps the original General
# formatting/style while its stored value remains a text string.
$updates = @{
    "D2" = "274.87"
    "E2" = "0.46%"
    "D3" = "27.22"
    "E3" = "2.13%"
    "D4" = "4.857"
    "E4" = "-0.29%"
    "D5" = "0.06393"
    "D6" = "6.949"
    "E6" = "1.14%"
    "E7" = "-1.08%"
    "D8" = "0.8758"
    "E8" = "0.64%"
    "D9" = "0.1511"
    "E9" = "3.77%"
    "D10" = "0.05118"
    "E10" = "0.63%"
    "D11" = "0.07587"
    "E11" = "2.84%"
    "D12" = "0.02961"
    "E12" = "-0.27%"
    "D13" = "0.08966"
    "E13" = "-0.79%"
    "D14" = "0.001565"
    "E14" = "-0.42%"
    "D15" = "0.0006398"
    "E15" = "1.90%"
    "D16" = "0.006188"
    "E16" = "2.75%"
    "E17" = "0.51%"
    "D18" = "3.308"
    "E18" = "-1.09%"
    "D19" = "2.283"
    "E19" = "0.49%"
    "E20" = "-1.09%"
    "E21" = "1.84%"
    "D22" = "3.907"
    "E22" = "-0.37%"
    "D23" = "0.04418"
    "E23" = "0.70%"
    "D25" = "0.001177"
    "E25" = "0.07%"
    "D26" = "0.003855"
    "E26" = "-9.66%"
    "D27" = "0.0001200"
    "E27" = "0.00%"
    "E28" = "14.64%"
    "D40" = "0.04125"
    "E40" = "2.46%"
    "D41" = "0.006789"
    "E41" = "1.34%"
    "D42" = "0.1171"
    "E42" = "0.37%"
    "D43" = "0.002110"
    "E43" = "0.48%"
    "D44" = "0.01147"
    "E44" = "-8.49%"
    "D45" = "0.00005186"
    "E45" = "-2.26%"
    "D46" = "1.679"
    "E46" = "-28.82%"
    "D47" = "0.01852"
    "E47" = "-7.38%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
